$wb = $excel.ActiveWorkbook

# --- Overview sheet (report for handoff) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-25 06:46:54"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Row 2 (b5670171-...): Status flips to "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
# Row 3 (bdad9598-...): Status + handoff datetime + error detail
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-25 06:46:49"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/169fc900c57033cd205f8fa1e454807d20c78cad/e2e/bdad9598-f144-4330-ba5a-155f7d09092c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b09449d4c8f23fe51fea1d8fbbac0fd038c991a/e2e/bdad9598-f144-4330-ba5a-155f7d09092c.md."
# Widen the Error Detail column to fit the new long text
$wsZhCn.Columns.Item(16).ColumnWidth = 39.15

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Row 3 (bdad9598-...): Status + handoff datetime + error detail
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-25 06:46:54"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/169fc900c57033cd205f8fa1e454807d20c78cad/e2e/bdad9598-f144-4330-ba5a-155f7d09092c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b09449d4c8f23fe51fea1d8fbbac0fd038c991a/e2e/bdad9598-f144-4330-ba5a-155f7d09092c.md."
# Widen the Error Detail column to fit the new long text
$wsDeDe.Columns.Item(16).ColumnWidth = 39.15
